$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 values (behavior profile recalculated) ---
$ws.Range("E2").Value = 50
$ws.Range("F2").Value = 4
$ws.Range("K2").Value = 490700
$ws.Range("N2").Value = 1298746.816666668

# --- Insert new row 3 (behavior id 2) ---
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 27
$ws.Range("D3").Value = 20
$ws.Range("E3").Value = 50
$ws.Range("F3").Value = 4
$ws.Range("G3").Value = 0.5
$ws.Range("H3").Value = 30
$ws.Range("I3").Value = "°C"
$ws.Range("J3").Value = 2
$ws.Range("K3").Value = 490700
$ws.Range("L3").Value = "Wh/person"
$ws.Range("M3").Value = 2
$ws.Range("N3").Value = 1368629.9833333339
$ws.Range("O3").Value = "Wh/person"
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 1

# --- Insert new row 4 (behavior id 3) ---
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 27
$ws.Range("D4").Value = 20
$ws.Range("E4").Value = 50
$ws.Range("F4").Value = 4
$ws.Range("G4").Value = 0.5
$ws.Range("H4").Value = 30
$ws.Range("I4").Value = "°C"
$ws.Range("J4").Value = 3
$ws.Range("K4").Value = 490700
$ws.Range("L4").Value = "Wh/person"
$ws.Range("M4").Value = 3
$ws.Range("N4").Value = 1423140.6500000004
$ws.Range("O4").Value = "Wh/person"
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 1

# --- View state matches the re-saved workbook (zoom + active selection) ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 132
[void]$ws.Range("F11").Select()
